# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (Exhibition) and "全部类型" (All types) worksheets, reflecting the
# freshly re-generated gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" — rows keyed by their F-column cell reference.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14891
$ws1.Range("F3").Value = 18577
$ws1.Range("F21").Value = 230
$ws1.Range("F22").Value = 7704
$ws1.Range("F26").Value = 1223
$ws1.Range("F28").Value = 5963
$ws1.Range("F34").Value = 5320

# Sheet "全部类型" — same events, different row numbers (no 演出 rows interleaved in the 展览-only sheet).
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14891
$ws4.Range("F3").Value = 18577
$ws4.Range("F22").Value = 230
$ws4.Range("F23").Value = 7704
$ws4.Range("F27").Value = 1223
$ws4.Range("F31").Value = 5963
$ws4.Range("F37").Value = 5320
